$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the newly-added day's data in row 50 (Date, Total Count,
# Session Timeout Errors, Errors Requiring Analysis)
$ws.Range("A50").Value = 46009
$ws.Range("B50").Value = 717
$ws.Range("C50").Value = 17
$ws.Range("D50").Value = 700

# Move the active selection down to the newly filled row, matching the
# saved sheet view state (activeCell A50, selection A50:D50)
$ws.Range("A50:D50").Select()
